# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for Chirimoya (Vega Modelo de Temuco)
# at sheet row 127, pushing the existing rows (127..194) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 127 (shifts 127:194 -> 128:195)
$ws.Rows(127).Insert()

# Populate the newly inserted row with the new observation
$ws.Range("A127").Value = 10
$ws.Range("B127").Value = "Vega Modelo de Temuco"
$ws.Range("C127").Value = "La Araucanía"
$ws.Range("D127").Value = 45161
$ws.Range("E127").Value = 9
$ws.Range("F127").Value = "Fruta"
$ws.Range("G127").Value = 100107
$ws.Range("H127").Value = "Otros"
$ws.Range("I127").Value = 100107002
$ws.Range("J127").Value = "Chirimoya"
$ws.Range("K127").Value = "Cultivar IV Región"
$ws.Range("L127").Value = "Primera"
$ws.Range("M127").Value = 90
$ws.Range("N127").Value = 3000
$ws.Range("O127").Value = 3000
$ws.Range("P127").Value = 3000
$ws.Range("Q127").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R127").Value = "Provincia del Elquí"
$ws.Range("S127").Value = 3000
$ws.Range("T127").Value = 1
